# Update Pergantian Ban & Service Kendaraan - 28/03/2024 15:45
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Zoom the sheet view to 110%
$excel.ActiveWindow.Zoom = 110

# Fill in row 22 with the new service entry
$ws.Range("B22").Value = "Fuel Filter"
$ws.Range("C22").Value = "Besar + Kecil"
$ws.Range("D22").Value = "2 pcs"
$ws.Range("E22").Value = 288107
$ws.Range("F22").Formula = "=10000+E22"
$ws.Range("G22").Value = "28/3/2024"
